$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column values that look like plain decimal numbers would silently
# be auto-converted to a Number by Excel (dropping trailing zeros / using
# scientific notation), losing fidelity with the source text values. Force text
# entry for those via NumberFormat "@", then ClearFormats() so the cell keeps no
# explicit style (matching the original unstyled cell) once the text is locked in.

$ws.Range("D2").Value = "29.394.86"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.848.61"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07598"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2929"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.850.51"
$ws.Range("E12").Value = "  -6.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001080"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +9.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6784"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.69"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.098.52"
$ws.Range("E17").Value = "  -7.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.177"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "29.409.12"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.480"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1396"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.343"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.463"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.299"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.104"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.028"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.842"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7092"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.588"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "1.239.23"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.417"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9058"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.77"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.172"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4021"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.030"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("E51").Value = "  -0.50%  "

# Rows 39 and 40 swapped coin identity (VeChain and MXToken traded ranking spots)
# with freshly scraped price/volume data.
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.776"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01801"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "
